$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "End Date" (Q2) was wrongly duplicating the "Expiry Date 2" (Z2) value
# ("26-6-2021") because of a copy/paste mistake upstream. Fix the duplicate
# by giving Q2 its correct reporting-month value instead.
# Force text (NumberFormat "@") while assigning so "Jun-2021" isn't
# auto-converted into a date serial, then restore General formatting to
# match the cell's original style.
$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = "Jun-2021"
$ws.Range("Q2").NumberFormat = "General"

# Leave the selection where the author last left it when saving.
$ws.Range("Q3").Select()
